$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1.Name = "Separate Testing Method"
$ws2.Name = "Cross Validation Method"

# --- Build out "Cross Validation Method" (formerly Sheet2) ---
# Column widths (approximate best-fit widths from the source workbook).
$ws2.Columns.Item(1).ColumnWidth = 25.5
$ws2.Columns.Item(2).ColumnWidth = 18.5
$ws2.Columns.Item(3).ColumnWidth = 10

# Cell values -- order matters so new shared-string entries land at the
# same indices as the target workbook.
$ws2.Range("A1").Value = "Model and Parameters"
$ws2.Range("A10").Value = "Note:"
$ws2.Range("A11").Value = "5-fold Cross Validation with averaged Accuracy, with New Data Split"

$ws2.Range("A2").Value = "MLPClassifier:`n  Layers (8-40, 8-40), 6 bins"
$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 34

$ws2.Range("B1").Value = "Best Parameter"
$ws2.Range("B2").Value = "Layers (33, 33)"

$ws2.Range("C1").Value = "Accuracy"
$ws2.Range("C1").NumberFormat = "0.0000000"
$ws2.Range("C2").NumberFormat = "0.0000000"
$ws2.Range("C2").Value = 0.99971585985737299

# --- Selection / active sheet bookkeeping ---
$ws1.Range("K6").Select()
$ws2.Activate()
$ws2.Range("A3").Select()
$excel.ActiveWindow.Zoom = 150
